$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4391
$ws.Range("I138").Value = 1649.6666
$ws.Range("J138").Value = 6218.5557
$ws.Range("K138").Value = 4948.9998
$ws.Range("L138").Value = 18655.6671
$ws.Range("M138").Value = 191.0002000000004
$ws.Range("N138").Value = -28935.6671

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 319755.72
$ws.Range("I32").Value = 3019.0513
$ws.Range("J32").Value = 2220175.8
$ws.Range("K32").Value = 3019.0513
$ws.Range("L32").Value = 2220175.8
$ws.Range("M32").Value = -2732.0513
$ws.Range("N32").Value = -2220749.8
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("N71").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 25020000
$ws.Range("I6").Value = 25020000
$ws.Range("K6").Value = 25020000
$ws.Range("M6").Value = -25019887
$ws.Range("H7").Value = 80.09999999999999
$ws.Range("I7").Value = 59.8
$ws.Range("J7").Value = 100.4
$ws.Range("K7").Value = 59.8
$ws.Range("L7").Value = 100.4
$ws.Range("M7").Value = 53.2
$ws.Range("N7").Value = -326.4
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("H25").Value = 15000
$ws.Range("J25").Value = 15000
$ws.Range("L25").Value = 15000
$ws.Range("N25").Value = -15348
$ws.Range("H31").Value = 1132512.4
$ws.Range("I31").Value = 2151940.8
$ws.Range("J31").Value = 3859.5715
$ws.Range("K31").Value = 2151940.8
$ws.Range("L31").Value = 3859.5715
$ws.Range("M31").Value = -2151645.8
$ws.Range("N31").Value = -4449.5715
$ws.Range("H34").Value = 1132512.4
$ws.Range("I34").Value = 2151940.8
$ws.Range("J34").Value = 3859.5715
$ws.Range("K34").Value = 2151940.8
$ws.Range("L34").Value = 3859.5715
$ws.Range("M34").Value = -2151738.8
$ws.Range("N34").Value = -4263.5715
$ws.Range("H41").Value = 3333
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("H51").Value = 19000
$ws.Range("J51").Value = 19000
$ws.Range("L51").Value = 19000
$ws.Range("N51").Value = -20472
$ws.Range("H59").Value = 11295.25
$ws.Range("I59").Value = 9000
$ws.Range("J59").Value = 12060.333
$ws.Range("K59").Value = 9000
$ws.Range("L59").Value = 12060.333
$ws.Range("M59").Value = -7855
$ws.Range("N59").Value = -14350.333
$ws.Range("H60").Value = 1090
$ws.Range("I60").Value = 1090
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 1090
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -579
$ws.Range("H61").Value = 19000
$ws.Range("J61").Value = 19000
$ws.Range("L61").Value = 19000
$ws.Range("N61").Value = -19696
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("H74").Value = 10142.5
$ws.Range("I74").Value = 10142.5
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 10142.5
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -9268.5
$ws.Range("H77").Value = 10142.5
$ws.Range("I77").Value = 10142.5
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 30427.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -26059.5
$ws.Range("H132").Value = 6253190
$ws.Range("I132").Value = 2680
$ws.Range("J132").Value = 16670706
$ws.Range("K132").Value = 8040
$ws.Range("L132").Value = 50012118
$ws.Range("M132").Value = -5510
$ws.Range("N132").Value = -50017178
$ws.Range("N17").ClearContents()
$ws.Range("N41").ClearContents()
$ws.Range("N50").ClearContents()
$ws.Range("N60").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("N74").ClearContents()
$ws.Range("N77").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2381.647
$ws.Range("J39").Value = 2381.647
$ws.Range("L39").Value = 7144.941
$ws.Range("N39").Value = -7732.941
$ws.Range("H74").Value = 11242
$ws.Range("J74").Value = 13999.25
$ws.Range("L74").Value = 41997.75
$ws.Range("N74").Value = -44119.75
$ws.Range("H77").Value = 11242
$ws.Range("J77").Value = 13999.25
$ws.Range("L77").Value = 125993.25
$ws.Range("N77").Value = -136601.25
$ws.Range("H113").Value = 829.2471
$ws.Range("I113").Value = 575.5
$ws.Range("J113").Value = 897.4179
$ws.Range("K113").Value = 1726.5
$ws.Range("L113").Value = 2692.2537
$ws.Range("M113").Value = 443.5
$ws.Range("N113").Value = -7032.2537
$ws.Range("H140").Value = 6699.3887
$ws.Range("I140").Value = 5756.4287
$ws.Range("J140").Value = 9999.75
$ws.Range("K140").Value = 17269.2861
$ws.Range("L140").Value = 29999.25
$ws.Range("M140").Value = -12089.2861
$ws.Range("N140").Value = -40359.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1512.6818
$ws.Range("I7").Value = 1293.2667
$ws.Range("K7").Value = 1293.2667
$ws.Range("M7").Value = -1181.2667
$ws.Range("H63").Value = 7375
$ws.Range("J63").Value = 7375
$ws.Range("L63").Value = 7375
$ws.Range("N63").Value = -8873
$ws.Range("H66").Value = 7375
$ws.Range("J66").Value = 7375
$ws.Range("L66").Value = 22125
$ws.Range("N66").Value = -29613
$ws.Range("H126").Value = 1512.6818
$ws.Range("I126").Value = 1293.2667
$ws.Range("K126").Value = 3879.800099999999
$ws.Range("M126").Value = -1409.800099999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 35717864
$ws.Range("I132").Value = 57693900
$ws.Range("K132").Value = 173081700
$ws.Range("M132").Value = -173079170
